$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").Value = '''27.251.95'
$ws.Range("E2").Value = '  -1.59%  '

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").Value = '''1.562.79'

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").Value = '  +0.13%  '

# Row 5: BNB -> BNB
$ws.Range("D5").Value = '''206.12'
$ws.Range("E5").Value = '  -0.49%  '

# Row 6: XRP -> XRP
$ws.Range("D6").Value = '''0.495'
$ws.Range("E6").Value = '  -1.37%  '

# Row 7: USDC -> USDC
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.15%  '

# Row 8: Solana -> Solana
$ws.Range("D8").Value = '''22.08'
$ws.Range("E8").Value = '  -0.36%  '

# Row 9: Cardano -> Cardano
$ws.Range("D9").Value = '''0.248'

# Row 10: Dogecoin -> Dogecoin
$ws.Range("D10").Value = '''0.0590'
$ws.Range("E10").Value = '  -0.17%  '

# Row 11: TRON -> TRON
$ws.Range("D11").Value = '''0.0865'
$ws.Range("E11").Value = '  -0.12%  '

# Row 12: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '''1.784.28'
$ws.Range("E12").Value = '  -1.40%  '

# Row 13: WrappedEther -> WrappedEther
$ws.Range("D13").Value = '''1.569.81'
$ws.Range("E13").Value = '  -0.49%  '

# Row 14: Polkadot -> Polkadot
$ws.Range("E14").Value = '  -1.80%  '

# Row 15: Polygon -> Polygon
$ws.Range("D15").Value = '''0.516'
$ws.Range("E15").Value = '  -2.39%  '

# Row 16: Litecoin -> Litecoin
$ws.Range("D16").Value = '''63.16'
$ws.Range("E16").Value = '  -0.58%  '

# Row 17: WrappedBTC -> WrappedBTC
$ws.Range("D17").Value = '''27.231.07'
$ws.Range("E17").Value = '  -1.45%  '

# Row 18: ShibaInu -> ShibaInu
$ws.Range("D18").Value = '''0.0₃0688'
$ws.Range("E18").Value = '  -0.96%  '

# Row 19: BitcoinCash -> BitcoinCash
$ws.Range("D19").Value = '''210.81'
$ws.Range("E19").Value = '  -4.02%  '

# Row 20: Chainlink -> Chainlink
$ws.Range("D20").Value = '''7.22'
$ws.Range("E20").Value = '  -1.17%  '

# Row 21: Dai -> Dai
$ws.Range("E21").Value = '  +0.08%  '

# Row 22: Uniswap -> Uniswap
$ws.Range("E22").Value = '  -0.80%  '

# Row 23: Avalanche -> Avalanche
$ws.Range("D23").Value = '''9.42'
$ws.Range("E23").Value = '  -1.34%  '

# Row 24: Toncoin -> Toncoin
$ws.Range("D24").Value = '''2.00'
$ws.Range("E24").Value = '  +1.73%  '

# Row 25: Monero -> Monero
$ws.Range("D25").Value = '''152.76'
$ws.Range("E25").Value = '  -0.57%  '

# Row 26: Cosmos -> Cosmos
$ws.Range("D26").Value = '''6.62'
$ws.Range("E26").Value = '  -3.47%  '

# Row 27: EthereumClassic -> EthereumClassic
$ws.Range("D27").Value = '''14.86'
$ws.Range("E27").Value = '  -1.58%  '

# Row 28: BinanceUSD -> BinanceUSD
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  +0.16%  '

# Row 29: Stellar -> Stellar
$ws.Range("E29").Value = '  -1.53%  '

# Row 30: PancakeSwap -> PancakeSwap
$ws.Range("E30").Value = '  -0.89%  '

# Row 31: Hedera -> Hedera
$ws.Range("E31").Value = '  -0.07%  '

# Row 32: Filecoin -> Filecoin
$ws.Range("D32").Value = '''3.16'
$ws.Range("E32").Value = '  -1.68%  '

# Row 33: Maker -> Maker
$ws.Range("D33").Value = '''1.383.17'
$ws.Range("E33").Value = '  +1.49%  '

# Row 34: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("E34").Value = '  +0.12%  '

# Row 35: LidoDAOToken -> LidoDAOToken
$ws.Range("E35").Value = '  +0.88%  '

# Row 36: HuobiToken -> HuobiToken
$ws.Range("D36").Value = '''2.30'
$ws.Range("E36").Value = '  -0.07%  '

# Row 37: TrustWalletToken -> TrustWalletToken
$ws.Range("D37").Value = '''0.940'
$ws.Range("E37").Value = '  -3.29%  '

# Row 38: VeChain -> VeChain
$ws.Range("E38").Value = '  -1.31%  '

# Row 39: ImmutableX -> ImmutableX
$ws.Range("D39").Value = '''0.522'
$ws.Range("E39").Value = '  -2.60%  '

# Row 40: ARBITRUM -> ARBITRUM
$ws.Range("D40").Value = '''0.813'
$ws.Range("E40").Value = '  -0.91%  '

# Row 41: PaxDollar -> PaxDollar
$ws.Range("E41").Value = '  +0.16%  '

# Row 42: WEMIXToken -> WEMIXToken
$ws.Range("E42").Value = '  +2.78%  '

# Row 43: MXToken -> RenderToken
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''1.81'
$ws.Range("E43").Value = '  +4.13%  '

# Row 44: RenderToken -> MXToken
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '''2.17'
$ws.Range("E44").Value = '  -0.19%  '

# Row 45: Aave -> Aave
$ws.Range("D45").Value = '''63.47'
$ws.Range("E45").Value = '  -0.43%  '

# Row 46: FraxShare -> FraxShare
$ws.Range("D46").Value = '''5.22'
$ws.Range("E46").Value = '  +0.40%  '

# Row 47: RocketPoolETH -> RocketPoolETH
$ws.Range("D47").Value = '''1.696.80'
$ws.Range("E47").Value = '  -1.34%  '

# Row 48: Quant -> Quant
$ws.Range("D48").Value = '''85.39'
$ws.Range("E48").Value = '  -2.97%  '

# Row 49: Cronos -> BabyDogeCoin
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.0₇0997'
$ws.Range("E49").Value = '  -1.03%  '

# Row 50: USDD -> Cronos
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0492'
$ws.Range("E50").Value = '  -0.99%  '

# Row 51: Algorand -> USDD
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '''1.01'
$ws.Range("E51").Value = '  +0.46%  '

